$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new data rows right before current row 400. This shifts the
# existing rows 400-407 down to become rows 405-412 (their content is left
# untouched), and the freed rows 400-404 are populated with new weekly
# price records below.
$ws.Range("A400:T404").EntireRow.Insert(-4121) | Out-Null

# New row 400
$ws.Range("A400").Value = 10
$ws.Range("B400").Value = "Vega Modelo de Temuco"
$ws.Range("C400").Value = "La Araucanía"
$ws.Range("D400").Value = 44448
$ws.Range("E400").Value = 9
$ws.Range("F400").Value = "Fruta"
$ws.Range("G400").Value = 100104
$ws.Range("H400").Value = "Frutos de pepita"
$ws.Range("I400").Value = 100104005
$ws.Range("J400").Value = "Pera"
$ws.Range("K400").Value = "Packham's Triumph"
$ws.Range("L400").Value = "Calibre 80"
$ws.Range("M400").Value = 500
$ws.Range("N400").Value = 19000
$ws.Range("O400").Value = 19000
$ws.Range("P400").Value = 19000
$ws.Range("Q400").Value = "$/caja 18 kilos embalada"
$ws.Range("R400").Value = "Región de O'Higgins"
$ws.Range("S400").Value = 1056
$ws.Range("T400").Value = 18

# New row 401
$ws.Range("A401").Value = 10
$ws.Range("B401").Value = "Vega Modelo de Temuco"
$ws.Range("C401").Value = "La Araucanía"
$ws.Range("D401").Value = 44448
$ws.Range("E401").Value = 9
$ws.Range("F401").Value = "Fruta"
$ws.Range("G401").Value = 100104
$ws.Range("H401").Value = "Frutos de pepita"
$ws.Range("I401").Value = 100104005
$ws.Range("J401").Value = "Pera"
$ws.Range("K401").Value = "Packham's Triumph"
$ws.Range("L401").Value = "Especial"
$ws.Range("M401").Value = 3
$ws.Range("N401").Value = 300000
$ws.Range("O401").Value = 300000
$ws.Range("P401").Value = 300000
$ws.Range("Q401").Value = "$/bins (450 kilos)"
$ws.Range("R401").Value = "Región de O'Higgins"
$ws.Range("S401").Value = 667
$ws.Range("T401").Value = 450

# New row 402
$ws.Range("A402").Value = 10
$ws.Range("B402").Value = "Vega Modelo de Temuco"
$ws.Range("C402").Value = "La Araucanía"
$ws.Range("D402").Value = 44448
$ws.Range("E402").Value = 9
$ws.Range("F402").Value = "Fruta"
$ws.Range("G402").Value = 100104
$ws.Range("H402").Value = "Frutos de pepita"
$ws.Range("I402").Value = 100104005
$ws.Range("J402").Value = "Pera"
$ws.Range("K402").Value = "Packham's Triumph"
$ws.Range("L402").Value = "Primera"
$ws.Range("M402").Value = 470
$ws.Range("N402").Value = 12000
$ws.Range("O402").Value = 13000
$ws.Range("P402").Value = 12468
$ws.Range("Q402").Value = "$/bandeja 18 kilos granel"
$ws.Range("R402").Value = "Región de O'Higgins"
$ws.Range("S402").Value = 693
$ws.Range("T402").Value = 18

# New row 403
$ws.Range("A403").Value = 10
$ws.Range("B403").Value = "Vega Modelo de Temuco"
$ws.Range("C403").Value = "La Araucanía"
$ws.Range("D403").Value = 44448
$ws.Range("E403").Value = 9
$ws.Range("F403").Value = "Fruta"
$ws.Range("G403").Value = 100104
$ws.Range("H403").Value = "Frutos de pepita"
$ws.Range("I403").Value = 100104005
$ws.Range("J403").Value = "Pera"
$ws.Range("K403").Value = "Packham's Triumph"
$ws.Range("L403").Value = "Primera"
$ws.Range("M403").Value = 8
$ws.Range("N403").Value = 240000
$ws.Range("O403").Value = 240000
$ws.Range("P403").Value = 240000
$ws.Range("Q403").Value = "$/bins (450 kilos)"
$ws.Range("R403").Value = "Región de O'Higgins"
$ws.Range("S403").Value = 533
$ws.Range("T403").Value = 450

# New row 404
$ws.Range("A404").Value = 10
$ws.Range("B404").Value = "Vega Modelo de Temuco"
$ws.Range("C404").Value = "La Araucanía"
$ws.Range("D404").Value = 44448
$ws.Range("E404").Value = 9
$ws.Range("F404").Value = "Fruta"
$ws.Range("G404").Value = 100104
$ws.Range("H404").Value = "Frutos de pepita"
$ws.Range("I404").Value = 100104005
$ws.Range("J404").Value = "Pera"
$ws.Range("K404").Value = "Packham's Triumph"
$ws.Range("L404").Value = "Segunda"
$ws.Range("M404").Value = 5
$ws.Range("N404").Value = 220000
$ws.Range("O404").Value = 220000
$ws.Range("P404").Value = 220000
$ws.Range("Q404").Value = "$/bins (450 kilos)"
$ws.Range("R404").Value = "Región de O'Higgins"
$ws.Range("S404").Value = 489
$ws.Range("T404").Value = 450
